$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 618.2
$ws.Range("I20").Value = 618.2
$ws.Range("K20").Value = 618.2
$ws.Range("M20").Value = -388.2

$ws.Range("H35").Value = 618.2
$ws.Range("I35").Value = 618.2
$ws.Range("K35").Value = 618.2
$ws.Range("M35").Value = -239.2

$ws.Range("H40").Value = 4823.2666
$ws.Range("I40").Value = 4700
$ws.Range("J40").Value = 4842.231
$ws.Range("K40").Value = 4700
$ws.Range("L40").Value = 4842.231
$ws.Range("M40").Value = -4525
$ws.Range("N40").Value = -5192.231

$ws.Range("H47").Value = 7000
$ws.Range("I47").Value = 7000
$ws.Range("K47").Value = 7000
$ws.Range("M47").Value = -6028

$ws.Range("H58").Value = 3500
$ws.Range("I58").Value = 250
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 750
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -600
$ws.Range("N58").Value = -30300

$ws.Range("H62").Value = 6867.5
$ws.Range("I62").Value = 5999
$ws.Range("K62").Value = 5999
$ws.Range("M62").Value = -5375

$ws.Range("H65").Value = 6867.5
$ws.Range("I65").Value = 5999
$ws.Range("K65").Value = 29995
$ws.Range("M65").Value = -26875

$ws.Range("H87").Value = 39999.332
$ws.Range("J87").Value = 39999.332
$ws.Range("L87").Value = 39999.332
$ws.Range("N87").Value = -42495.332

$ws.Range("H90").Value = 39999.332
$ws.Range("J90").Value = 39999.332
$ws.Range("L90").Value = 119997.996
$ws.Range("N90").Value = -132477.996

$ws.Range("H107").Value = 1070.6765
$ws.Range("J107").Value = 2382.6667
$ws.Range("L107").Value = 2382.6667
$ws.Range("N107").Value = -6222.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 348
$ws.Range("I22").Value = 334.875
$ws.Range("J22").Value = 400.5
$ws.Range("K22").Value = 334.875
$ws.Range("L22").Value = 400.5
$ws.Range("M22").Value = -161.875
$ws.Range("N22").Value = -746.5

$ws.Range("H49").Value = 2000
$ws.Range("J49").Value = 2000
$ws.Range("L49").Value = 2000
$ws.Range("N49").Value = -2478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7122
$ws.Range("I16").Value = 5983
$ws.Range("J16").Value = 9400
$ws.Range("K16").Value = 5983
$ws.Range("L16").Value = 9400
$ws.Range("M16").Value = -5696
$ws.Range("N16").Value = -9974

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H41").Value = 17007.691
$ws.Range("J41").Value = 18081.818
$ws.Range("L41").Value = 18081.818
$ws.Range("N41").Value = -18937.818

$ws.Range("H113").Value = 7122
$ws.Range("I113").Value = 5983
$ws.Range("J113").Value = 9400
$ws.Range("K113").Value = 5983
$ws.Range("L113").Value = 9400
$ws.Range("M113").Value = -3813
$ws.Range("N113").Value = -13740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 400
$ws.Range("J29").Value = 400
$ws.Range("L29").Value = 1200
$ws.Range("N29").Value = -1754

$ws.Range("H41").Value = 3000
$ws.Range("J41").Value = 3000
$ws.Range("L41").Value = 9000
$ws.Range("N41").Value = -9676

$ws.Range("H50").Value = 398.4
$ws.Range("I50").Value = 518
$ws.Range("J50").Value = 119.333336
$ws.Range("K50").Value = 1554
$ws.Range("L50").Value = 358.000008
$ws.Range("M50").Value = -1073
$ws.Range("N50").Value = -1320.000008

$ws.Range("H53").Value = 398.4
$ws.Range("I53").Value = 518
$ws.Range("J53").Value = 119.333336
$ws.Range("K53").Value = 1554
$ws.Range("L53").Value = 358.000008
$ws.Range("M53").Value = -1073
$ws.Range("N53").Value = -1320.000008

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H113").Value = 1839
$ws.Range("J113").Value = 2399
$ws.Range("L113").Value = 7197
$ws.Range("N113").Value = -11537

$ws.Range("H129").Value = 3943
$ws.Range("J129").Value = 3943
$ws.Range("L129").Value = 11829
$ws.Range("N129").Value = -21829

$ws.Range("H131").Value = 14728.177
$ws.Range("I131").Value = 30522
$ws.Range("J131").Value = 3672.5
$ws.Range("K131").Value = 91566
$ws.Range("L131").Value = 11017.5
$ws.Range("M131").Value = -86526
$ws.Range("N131").Value = -21097.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 15000
$ws.Range("J39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -16064

$ws.Range("H113").Value = 1109.1666
$ws.Range("I113").Value = 1109.1666
$ws.Range("K113").Value = 1109.1666
$ws.Range("M113").Value = 1060.8334

$ws.Range("H123").Value = 4029999.5
$ws.Range("J123").Value = 4029999.5
$ws.Range("L123").Value = 4029999.5
$ws.Range("N123").Value = -4034899.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1666.3334
$ws.Range("I22").Value = 498.5
$ws.Range("K22").Value = 498.5
$ws.Range("M22").Value = -203.5

$ws.Range("H27").Value = 1666.3334
$ws.Range("I27").Value = 498.5
$ws.Range("K27").Value = 498.5
$ws.Range("M27").Value = -391.5

$ws.Range("H40").Value = 4006.158
$ws.Range("I40").Value = 4210.5
$ws.Range("J40").Value = 2916.3333
$ws.Range("K40").Value = 4210.5
$ws.Range("L40").Value = 2916.3333
$ws.Range("M40").Value = -4074.5
$ws.Range("N40").Value = -3188.3333

$ws.Range("H46").Value = 4436.375
$ws.Range("J46").Value = 4570.143
$ws.Range("L46").Value = 4570.143
$ws.Range("N46").Value = -4946.143

$ws.Range("H100").Value = 6211.7144
$ws.Range("I100").Value = 6699
$ws.Range("J100").Value = 4993.5
$ws.Range("K100").Value = 6699
$ws.Range("L100").Value = 4993.5
$ws.Range("M100").Value = -6158
$ws.Range("N100").Value = -6075.5

$ws.Range("H122").Value = 4315.8335
$ws.Range("I122").Value = 4274
$ws.Range("K122").Value = 12822
$ws.Range("M122").Value = -10372

$ws.Range("H136").Value = 6416.7144
$ws.Range("J136").Value = 7582
$ws.Range("L136").Value = 22746
$ws.Range("N136").Value = -27846
